# Auto-generated edit script: update '想去人数' (interest count) values
# in sheets 展览, 演出, and 全部类型 to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 70
$ws.Range("F6").Value = 2720
$ws.Range("F8").Value = 1601
$ws.Range("F9").Value = 7361
$ws.Range("F11").Value = 7529
$ws.Range("F14").Value = 5980
$ws.Range("F15").Value = 3216
$ws.Range("F16").Value = 3583
$ws.Range("F17").Value = 4
$ws.Range("F18").Value = 3
$ws.Range("F20").Value = 431
$ws.Range("F21").Value = 270
$ws.Range("F22").Value = 267
$ws.Range("F23").Value = 2066
$ws.Range("F26").Value = 919
$ws.Range("F28").Value = 943
$ws.Range("F30").Value = 2570
$ws.Range("F31").Value = 1392
$ws.Range("F32").Value = 3135
$ws.Range("F33").Value = 134
$ws.Range("F35").Value = 227
$ws.Range("F36").Value = 888
$ws.Range("F37").Value = 462
$ws.Range("F38").Value = 1207
$ws.Range("F41").Value = 572

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 51
$ws.Range("F6").Value = 41
$ws.Range("F9").Value = 391
$ws.Range("F13").Value = 10

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 70
$ws.Range("F7").Value = 51
$ws.Range("F10").Value = 2720
$ws.Range("F11").Value = 1601
$ws.Range("F14").Value = 7361
$ws.Range("F16").Value = 7529
$ws.Range("F18").Value = 5980
$ws.Range("F19").Value = 3216
$ws.Range("F20").Value = 3583
$ws.Range("F21").Value = 4
$ws.Range("F22").Value = 3
$ws.Range("F24").Value = 431
$ws.Range("F26").Value = 270
$ws.Range("F29").Value = 267
$ws.Range("F30").Value = 2066
$ws.Range("F31").Value = 10
$ws.Range("F36").Value = 943
$ws.Range("F38").Value = 2570
$ws.Range("F39").Value = 1393
$ws.Range("F41").Value = 3135
$ws.Range("F42").Value = 134
$ws.Range("F45").Value = 462
$ws.Range("F46").Value = 1207
$ws.Range("F49").Value = 572
